$d = $word.ActiveDocument

# --- 1. Tag the three image-only paragraphs as "no proofing" -----------------
# Word marks runs that hold an inline picture with <w:rPr><w:noProof/></w:rPr>.
# Walk every paragraph and flip NoProofing on for any run whose range contains
# an inline picture (there are three such paragraphs in this journal).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.Font.NoProofing = -1
    }
}

# --- 2. Relocate the "_GoBack" bookmark --------------------------------------
# It currently sits at the end of the paragraph "... Since I am the admin I
# see this:". Word re-drops this bookmark at the point of the most recent
# edit, which after this change is the trailing empty Heading-1 paragraph at
# the very end of the document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 3. Add the new sentence after "I merged the pull request." -------------
# Find that paragraph via its text so the script doesn't depend on a brittle
# paragraph index.
$mergedPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "I merged the pull request.") {
        $mergedPara = $p
    }
}

$r = $mergedPara.Range
$r.MoveEnd(1, -1)      # exclude the paragraph mark
$r.Collapse(0)         # collapse to the end of the existing sentence

# Dropping a (temporary) bookmark at the insertion point keeps the new text
# in its own run instead of being coalesced into the identically-formatted
# run that precedes it; the bookmark itself is removed again immediately
# after the text is in place.
$d.Bookmarks.Add("ZZZTempSplitMark", $r)
$newText = " There were some problems with that, do I am trying again. I" + [char]0x2019 + "m not certain what went wrong."
$r.InsertAfter($newText)
$d.Bookmarks("ZZZTempSplitMark").Delete()

# Now that the edit has been made, re-add "_GoBack" at the very last
# paragraph of the document, which is where Word leaves it after this kind
# of change.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
